# Update the workbook's "build version" string everywhere it appears.
#
# Old version string: "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
# New version string: "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

# "About" sheet: the version appears in the headline (A2) and inside the
# recommended-citation text (A6).
$wsAbout = $wb.Worksheets.Item("About")

$a2 = $wsAbout.Range("A2").Value2
$wsAbout.Range("A2").Value2 = $a2.Replace($oldVersion, $newVersion)

$a6 = $wsAbout.Range("A6").Value2
$wsAbout.Range("A6").Value2 = $a6.Replace($oldVersion, $newVersion)

# "Boundaries and methane sources" sheet: column S ("build_version") carries
# the same version string on every data row.
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

$lastRow = $wsData.Cells.Item($wsData.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = $wsData.UsedRange.Rows.Count
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $wsData.Cells.Item($row, 19) # column S
    $current = $cell.Value2
    if ($current -ne $null -and $current -is [string] -and $current.Contains($oldVersion)) {
        $cell.Value2 = $current.Replace($oldVersion, $newVersion)
    }
}
